$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (new report week / issue number) ---
$ws.Range("A8").Value = "Volume 30   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# --- Weekly crime-statistics table refresh (rows 14-30) ---
# Reference cells already bearing the target style+shared-string combos;
# used via Copy/PasteSpecial to reproduce the exact text-placeholder styling
# Excel applies to suppressed ("0") and not-applicable ("***.*") figures.
$srcZero = $ws.Range("C14")    # style: placeholder showing "0"
$srcNA   = $ws.Range("E14")    # style: placeholder showing "***.*"
$srcNum  = $ws.Range("J14")    # style: plain right-aligned integer

$srcZero.Copy()
$ws.Range("G14").PasteSpecial(-4163)
$srcZero.Copy()
$ws.Range("G14").PasteSpecial(-4122)
$srcNA.Copy()
$ws.Range("H14").PasteSpecial(-4163)
$srcNA.Copy()
$ws.Range("H14").PasteSpecial(-4122)
$srcZero.Copy()
$ws.Range("C15").PasteSpecial(-4163)
$srcZero.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 50
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -11.111111111111
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 167
$ws.Range("J16").Value = 125
$ws.Range("K16").Value = 33.6
$ws.Range("L16").Value = 114.102564102564
$ws.Range("M16").Value = 65.346534653465
$ws.Range("N16").Value = -69.964028776978
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 22.580645161290
$ws.Range("I17").Value = 232
$ws.Range("J17").Value = 135
$ws.Range("K17").Value = 71.851851851851
$ws.Range("L17").Value = 129.70297029703
$ws.Range("M17").Value = 136.734693877551
$ws.Range("N17").Value = 35.672514619883
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 31
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -3.125
$ws.Range("I18").Value = 279
$ws.Range("J18").Value = 221
$ws.Range("K18").Value = 26.244343891402
$ws.Range("L18").Value = 78.846153846153
$ws.Range("M18").Value = 27.397260273972
$ws.Range("N18").Value = -75.78125
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 52.631578947368
$ws.Range("F19").Value = 96
$ws.Range("G19").Value = 96
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 679
$ws.Range("J19").Value = 721
$ws.Range("K19").Value = -5.825242718446
$ws.Range("L19").Value = 115.555555555556
$ws.Range("M19").Value = 133.333333333333
$ws.Range("N19").Value = 5.763239875389
$ws.Range("C20").Value = 10
$ws.Range("E20").Value = 233.333333333333
$ws.Range("F20").Value = 32
$ws.Range("H20").Value = 166.666666666667
$ws.Range("I20").Value = 240
$ws.Range("J20").Value = 103
$ws.Range("K20").Value = 133.009708737864
$ws.Range("L20").Value = 128.571428571429
$ws.Range("M20").Value = 77.777777777777
$ws.Range("N20").Value = -87.355110642781
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 47.727272727272
$ws.Range("F21").Value = 224
$ws.Range("G21").Value = 203
$ws.Range("H21").Value = 10.344827586206
$ws.Range("I21").Value = 1615
$ws.Range("J21").Value = 1323
$ws.Range("K21").Value = 22.071050642479
$ws.Range("L21").Value = 110.01300390117
$ws.Range("M21").Value = 88.668224299065
$ws.Range("N21").Value = -63.609734114466
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 19
$ws.Range("K22").Value = 18.75
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = 533.333333333333
$srcZero.Copy()
$ws.Range("G23").PasteSpecial(-4163)
$srcZero.Copy()
$ws.Range("G23").PasteSpecial(-4122)
$srcNA.Copy()
$ws.Range("H23").PasteSpecial(-4163)
$srcNA.Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("L23").Value = 12.5
$ws.Range("M23").Value = 28.571428571428
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = -10.714285714285
$ws.Range("F24").Value = 214
$ws.Range("G24").Value = 231
$ws.Range("H24").Value = -7.359307359307
$ws.Range("I24").Value = 1343
$ws.Range("J24").Value = 1369
$ws.Range("K24").Value = -1.899196493791
$ws.Range("L24").Value = 71.301020408163
$ws.Range("M24").Value = 84.986225895316
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 144.444444444444
$ws.Range("F25").Value = 67
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 31.372549019607
$ws.Range("I25").Value = 402
$ws.Range("J25").Value = 327
$ws.Range("K25").Value = 22.935779816513
$ws.Range("L25").Value = 89.622641509434
$ws.Range("M25").Value = 16.860465116279
$srcZero.Copy()
$ws.Range("C26").PasteSpecial(-4163)
$srcZero.Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = 66.666666666666
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 55
$ws.Range("J27").Value = 45
$ws.Range("K27").Value = 22.222222222222
$ws.Range("L27").Value = 103.703703703704
$srcNum.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$srcZero.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$srcZero.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$srcNA.Copy()
$ws.Range("E28").PasteSpecial(-4163)
$srcNA.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 5
$ws.Range("K28").Value = 25
$ws.Range("L28").Value = 400
$ws.Range("M28").Value = 400
$ws.Range("N28").Value = -58.333333333333
$srcNum.Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$srcZero.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$srcZero.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$srcNA.Copy()
$ws.Range("E29").PasteSpecial(-4163)
$srcNA.Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 25
$ws.Range("L29").Value = 400
$ws.Range("M29").Value = 400
$ws.Range("N29").Value = -50
$ws.Range("F30").Value = 1

$excel.CutCopyMode = 0
